$d = $word.ActiveDocument

# 1) Merge the three title runs ("Mul" + "tiple Choice Questions Chapter 3" + ":")
#    into a single run with the same text.
$d.Content.Find.Execute("Multiple Choice Questions Chapter 3:", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Multiple Choice Questions Chapter 3:", 2) | Out-Null

# 2) Remove the old _GoBack bookmark from its original location
#    (right after "Obeys the requirements of the algorithm.")
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 3) Insert a new, empty paragraph directly after the existing blank paragraph
#    that follows the title, and place the (now-empty) _GoBack bookmark in it.
$p2 = $d.Paragraphs(2)
$p2.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs(3)

# Use a temporary placeholder character so the bookmark range resolves
# unambiguously to the interior of the new paragraph, then remove the
# placeholder, leaving a clean <w:p><w:bookmarkStart/><w:bookmarkEnd/></w:p>.
$insPt = $d.Range($newPara.Range.Start, $newPara.Range.Start)
$insPt.InsertAfter("X")
$newPara2 = $d.Paragraphs(3)
$bmRange = $d.Range($newPara2.Range.Start, $newPara2.Range.Start + 1)
$d.Bookmarks.Add("_GoBack", $bmRange)
$delRange = $d.Range($newPara2.Range.Start, $newPara2.Range.Start + 1)
$delRange.Delete()

# 4) Delete the entire first question ("The strength of a Hash function...")
#    together with its answer/reason paragraphs and the trailing blank line,
#    from "The strength of a Hash function..." through the blank paragraph
#    right before "In RSA, bits in 'e' and 'd' ...".
$startRange = $d.Range(0, 0)
$startRange.Find.Execute("The strength of a Hash function against brute force attack depends on") | Out-Null
$endRange = $d.Range(0, 0)
$endRange.Find.Execute("In RSA, bits in") | Out-Null

$deleteRange = $d.Range($startRange.Start, $endRange.Start)
$deleteRange.Delete()

# 5) Merge "A" + "nswer - c) Increases security of the algorithm." into one run.
$d.Content.Find.Execute("Answer – c) Increases security of the algorithm.", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Answer – c) Increases security of the algorithm.", 2) | Out-Null

# 6) Merge "Reason - " + "When the values of bits in 'e' and 'd' ..." into one run.
$d.Content.Find.Execute("Reason – When the values of bits in", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Reason – When the values of bits in", 2) | Out-Null
